$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B17").Value = "test"
